$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the two rows that were dropped from the dataset ---
# Row 26 = "RM 232", row 28 = "SC 92" (in the original layout).
# Deleting row 26 first shifts "SC 92" up to row 27, so delete row 27 next.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Step 2: fix up the "missing value" markers in column B ---
# After the row deletions, "SC 120" (row 30) gets its previously-missing B value filled in,
# while "SC 193" (row 32) now has its B value treated as missing instead.
$ws.Cells.Item(30, 2).Value = -19.7
$ws.Cells.Item(32, 2).Value = $null

# --- Step 3: update column F (the computed/imputed error column) ---
# Some rows gain a newly-computed value, others have their value cleared out.
$ws.Cells.Item(2, 6).Value = 18.03
$ws.Cells.Item(6, 6).Value = $null
$ws.Cells.Item(12, 6).Value = 17.45
$ws.Cells.Item(14, 6).Value = $null
$ws.Cells.Item(20, 6).Value = 17.73
$ws.Cells.Item(21, 6).Value = 16.58
$ws.Cells.Item(22, 6).Value = $null
$ws.Cells.Item(23, 6).Value = $null
$ws.Cells.Item(31, 6).Value = 17.18
$ws.Cells.Item(33, 6).Value = 17.53
